# "added version that takes in all 3 datatypes"
# Update two data values on Sheet1 and move the active selection to A3
# (also clears the old scrolled-down view/E28 selection in the process).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Data edits: A2 2 -> 18, A3 3 -> 8
$ws.Range("A2").Value = 18
$ws.Range("A3").Value = 8

# Restore the view to the top of the sheet and select A3 (was scrolled to
# A10/E28 before).
$ws.Activate()
$ws.Range("A3").Select()
